# Updated cryptos list on Mon Apr 15 20:51:29 UTC 2024 with GitHub Actions
# Refresh price / 1h-volume-change figures scraped from coinranking.com,
# plus two pairs of rows whose relative ranking swapped (Toncoin/Dogecoin
# at rows 10-11, Maker/Hedera at rows 39-40).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimals (e.g. "554.74") that Excel would
# otherwise auto-convert to numbers; force those specific cells to stay
# text (NumberFormat "@") so they match the original text-cell storage.
# Values with two dots (e.g. "63.256.34") or other non-numeric-looking
# text never get auto-converted, so no NumberFormat change is needed there.

$ws.Range('D2').Value = '63.256.34'
$ws.Range('E2').Value = '  -0.92%  '
$ws.Range('D3').Value = '3.088.95'
$ws.Range('E3').Value = '  +0.75%  '
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '554.74'
$ws.Range('E5').Value = '  +0.82%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '137.08'
$ws.Range('E6').Value = '  -3.00%  '
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '3.078.93'
$ws.Range('E8').Value = '  +0.62%  '
$ws.Range('E9').Value = '  +1.78%  '
$ws.Range('B10').Value = 'Toncoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.61'
$ws.Range('E10').Value = '  +2.76%  '
$ws.Range('B11').Value = 'Dogecoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.159'
$ws.Range('E11').Value = '  +6.14%  '
$ws.Range('E12').Value = '  +1.83%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '35.00'
$ws.Range('E13').Value = '  -0.69%  '
$ws.Range('E14').Value = '  +1.47%  '
$ws.Range('D15').Value = '3.582.46'
$ws.Range('E15').Value = '  +0.39%  '
$ws.Range('D16').Value = '63.248.59'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('E17').Value = '  +0.20%  '
$ws.Range('D18').Value = '3.089.84'
$ws.Range('E18').Value = '  +0.47%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '500.15'
$ws.Range('E19').Value = '  +3.29%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.66'
$ws.Range('E20').Value = '  +1.92%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.52'
$ws.Range('E21').Value = '  +0.33%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.706'
$ws.Range('E22').Value = '  +4.55%  '
$ws.Range('E23').Value = '  +1.81%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '78.01'
$ws.Range('E24').Value = '  +1.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '12.30'
$ws.Range('E25').Value = '  -0.04%  '
$ws.Range('E26').Value = '  +0.09%  '
$ws.Range('E27').Value = '  +2.32%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.16'
$ws.Range('E28').Value = '  -0.14%  '
$ws.Range('E29').Value = '  -1.25%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('E30').Value = '  -0.12%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '26.27'
$ws.Range('E31').Value = '  +2.99%  '
$ws.Range('E32').Value = '  -3.46%  '
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '59.14'
$ws.Range('E34').Value = '  +13.94%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '533.83'
$ws.Range('E35').Value = '  -7.01%  '
$ws.Range('E36').Value = '  +1.32%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.15'
$ws.Range('E37').Value = '  -2.44%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.0413'
$ws.Range('E38').Value = '  +3.88%  '
$ws.Range('B39').Value = 'Maker'
$ws.Range('C39').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D39').Value = '3.057.59'
$ws.Range('E39').Value = '  +2.34%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0791'
$ws.Range('E40').Value = '  +1.28%  '
$ws.Range('E41').Value = '  +2.84%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.07'
$ws.Range('E42').Value = '  -0.23%  '
$ws.Range('E43').Value = '  -5.21%  '
$ws.Range('E44').Value = '  +4.79%  '
$ws.Range('E45').Value = '  +0.00%  '
$ws.Range('E46').Value = '  -0.38%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '120.30'
$ws.Range('E47').Value = '  +1.83%  '
$ws.Range('E48').Value = '  +0.27%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '23.77'
$ws.Range('E49').Value = '  -4.01%  '
$ws.Range('D50').Value = '0.0₃0495'
$ws.Range('E50').Value = '  -4.14%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.34'
$ws.Range('E51').Value = '  +59.77%  '
